# Apply the changes described by the diff to the workbook.

$wb = $excel.ActiveWorkbook

# --- Worksheet "RequestedWorkDay" (sheet1.xml) -----------------------------
# It is no longer the tab-selected sheet; its remembered selection moves to D49.
$wsRWD = $wb.Worksheets.Item("RequestedWorkDay")
$wsRWD.Activate()
$wsRWD.Range("D49").Select()

# --- Worksheet "ListOfPA" (sheet3.xml) -------------------------------------
# Add the 16 PA names/shift-counts below the existing header row.
$wsPA = $wb.Worksheets.Item("ListOfPA")

$names = @(
    "Jun (A1)",
    "Jay (A2)",
    "Colleen (A3)",
    "Vivian (A4)",
    "Ralp (A5)",
    "Jimy (A6)",
    "Gavin (A7)",
    "Maggie (A8)",
    "JayZ (A9)",
    "Parker (A10)",
    "Singh (A11)",
    "Jupiter (A12)",
    "Kasi (A13)",
    "Asah (A14)",
    "Tyler (A15)",
    "Mike (A16)"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $wsPA.Cells.Item($row, 1).Value = $names[$i]
    $wsPA.Cells.Item($row, 2).Value = 12
}

# ListOfPA becomes the active/tab-selected sheet, zoomed to 175%, with B17
# selected (this also drives workbook.xml's bookViews/activeTab).
$wsPA.Activate()
$wsPA.Range("B17").Select()
$excel.ActiveWindow.Zoom = 175
